$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a plain text value without Excel auto-converting it to a number/date.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Cells whose new text will not be misread as a number (percentages, multi-dot prices, names, links) ---
$ws.Range("D2").Value = '62.850.26'
$ws.Range("E2").Value = '  +2.39%  '
$ws.Range("D3").Value = '3.472.17'
$ws.Range("E3").Value = '  +2.63%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("E6").Value = '  +4.67%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +1.94%  '
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("E10").Value = '  +2.33%  '
$ws.Range("E11").Value = '  +4.04%  '
$ws.Range("D12").Value = '4.063.89'
$ws.Range("E12").Value = '  +2.59%  '
$ws.Range("E13").Value = '  +5.39%  '
$ws.Range("E14").Value = '  +2.60%  '
$ws.Range("D15").Value = '3.471.56'
$ws.Range("E15").Value = '  +1.74%  '
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '62.865.41'
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("E18").Value = '  +4.46%  '
$ws.Range("E19").Value = '  +5.78%  '
$ws.Range("E20").Value = '  +2.89%  '
$ws.Range("E21").Value = '  -0.38%  '
$ws.Range("E22").Value = '  +2.58%  '
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = '3.605.57'
$ws.Range("E25").Value = '  +2.32%  '
$ws.Range("E26").Value = '  +3.57%  '
$ws.Range("E27").Value = '  -8.63%  '
$ws.Range("E28").Value = '  +4.75%  '
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").Value = '  +2.52%  '
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("E33").Value = '  -0.44%  '
$ws.Range("E34").Value = '  +2.42%  '
$ws.Range("E35").Value = '  +6.58%  '
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("E37").Value = '  +8.45%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E38").Value = '  +22.02%  '
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E39").Value = '  +1.38%  '
$ws.Range("D40").Value = '3.509.80'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("E42").Value = '  +2.92%  '
$ws.Range("E43").Value = '  +1.51%  '
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("E45").Value = '  +4.85%  '
$ws.Range("E46").Value = '  +4.45%  '
$ws.Range("D47").Value = '2.605.99'
$ws.Range("E47").Value = '  +6.07%  '
$ws.Range("E48").Value = '  +3.94%  '
$ws.Range("E49").Value = '  +1.47%  '
$ws.Range("E50").Value = '  +11.13%  '
$ws.Range("E51").Value = '  -0.03%  '

# --- Price cells whose new text looks like a plain number; force them to stay text ---
Set-TextValue "D4" '0.999'
Set-TextValue "D5" '580.19'
Set-TextValue "D6" '147.29'
Set-TextValue "D13" '29.67'
Set-TextValue "D14" '0.129'
Set-TextValue "D18" '6.37'
Set-TextValue "D19" '14.36'
Set-TextValue "D20" '9.23'
Set-TextValue "D21" '388.61'
Set-TextValue "D22" '0.562'
Set-TextValue "D23" '74.87'
Set-TextValue "D26" '0.0000116'
Set-TextValue "D27" '0.179'
Set-TextValue "D28" '7.63'
Set-TextValue "D29" '0.998'
Set-TextValue "D30" '8.16'
Set-TextValue "D31" '2.14'
Set-TextValue "D34" '23.74'
Set-TextValue "D35" '5.32'
Set-TextValue "D36" '7.08'
Set-TextValue "D37" '1.59'
Set-TextValue "D38" '31.50'
Set-TextValue "D39" '170.76'
Set-TextValue "D41" '0.0767'
Set-TextValue "D42" '0.801'
Set-TextValue "D44" '42.27'
Set-TextValue "D46" '1.20'
Set-TextValue "D48" '23.39'
Set-TextValue "D49" '6.74'
Set-TextValue "D50" '2.22'

Write-Host "Applied cryptos list update"
